$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 - new entry: topic learned + date
$ws.Range("C31").Value = "Temas aprendidos: Poo, Prototypes, clase, herencia y Try Catch"
$ws.Range("D29").Copy()
$ws.Range("D31").PasteSpecial(-4122)  # xlPasteFormats - keep the existing date style (s=7)
$ws.Range("D31").Value = 44747

# Row 33 - new entry: topic learned + date
$ws.Range("C33").Value = "Temas aprendidos:  Async Await, Fetch Api , introducción a Php"
$ws.Range("D29").Copy()
$ws.Range("D33").PasteSpecial(-4122)  # xlPasteFormats - keep the existing date style (s=7)
$ws.Range("D33").Value = 44748

# Update the view: scroll/selection moved
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 2
$ws.Range("C29:C30").Select()
